# Update TPM-derived values in the Plg-Itgb1 NATMI output sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ligand average/total expression values (same ligand -> same for every row)
$ws.Range("G2:G6").Value = 0.0345785
$ws.Range("H2:H6").Value = 0.069157

# Row 2 (Target cluster: ECs)
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 4.216112901511499
$ws.Range("R2").Value = 25.29667740906899
$ws.Range("S2").Value = 0.2282232151508951
$ws.Range("T2").Value = 0.2419720431319445

# Row 3 (Target cluster: FAPs) - M3/N3 unchanged
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 5.114655659904999
$ws.Range("R3").Value = 30.68793395942999
$ws.Range("S3").Value = 0.2768624053389947
$ws.Range("T3").Value = 0.2935413991166814

# Row 4 (Target cluster: Inflammatory-Mac)
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 2.887476582092666
$ws.Range("R4").Value = 17.324859492556
$ws.Range("S4").Value = 0.1563025480180701
$ws.Range("T4").Value = 0.1657186665504434

# Row 5 (Target cluster: MuSCs)
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 3.14901083055625
$ws.Range("R5").Value = 12.596043322225
$ws.Range("S5").Value = 0.1704597085236707
$ws.Range("T5").Value = 0.1204857969594293

# Row 6 (Target cluster: Resolving-Mac)
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 3.106381331955
$ws.Range("R6").Value = 18.63828799173
$ws.Range("S6").Value = 0.1681521229683693
$ws.Range("T6").Value = 0.1782820942415013
